$wb = $excel.ActiveWorkbook

# départements (Worksheets.Item(1))
$ws = $wb.Worksheets.Item(1)
$ws.Range("C97").Value = 70.58823529411765
$ws.Range("D97").Value = 36
$ws.Range("C99").Value = 15.11627906976744
$ws.Range("D99").Value = 13
$ws.Range("C104").Value = 4.6875
$ws.Range("D104").Value = 6
$ws.Range("C121").Value = 15
$ws.Range("D121").Value = 9
$ws.Range("C123").Value = 3.333333333333333
$ws.Range("D123").Value = 3
$ws.Range("E123").Value = 90
$ws.Range("C124").Value = 8.108108108108109
$ws.Range("D124").Value = 6
$ws.Range("E124").Value = 74
$ws.Range("C125").Value = 6.349206349206349
$ws.Range("D125").Value = 4
$ws.Range("C130").Value = 6.432748538011696
$ws.Range("D130").Value = 11
$ws.Range("C132").Value = 4.25531914893617
$ws.Range("D132").Value = 8
$ws.Range("C136").Value = 18.29268292682927
$ws.Range("D136").Value = 15
$ws.Range("C161").Value = 7.865168539325842
$ws.Range("D161").Value = 14
$ws.Range("C171").Value = 10.58823529411765
$ws.Range("D171").Value = 9
$ws.Range("C173").Value = 11.23595505617977
$ws.Range("D173").Value = 10
$ws.Range("C175").Value = 6.504065040650407
$ws.Range("D175").Value = 8
$ws.Range("C177").Value = 9.558823529411764
$ws.Range("D177").Value = 13
$ws.Range("C187").Value = 10.14492753623188
$ws.Range("D187").Value = 7
$ws.Range("C190").Value = 7.086614173228346
$ws.Range("D190").Value = 9
$ws.Range("C193").Value = 5.063291139240507
$ws.Range("D193").Value = 4
$ws.Range("C201").Value = 4.929577464788732
$ws.Range("D201").Value = 7
$ws.Range("C213").Value = 7.142857142857142
$ws.Range("D213").Value = 4
$ws.Range("C220").Value = 6.25
$ws.Range("D220").Value = 6
$ws.Range("E220").Value = 96
$ws.Range("C227").Value = 5.263157894736842
$ws.Range("D227").Value = 9
$ws.Range("C247").Value = 3.370786516853932
$ws.Range("D247").Value = 3
$ws.Range("C263").Value = 6.569343065693431
$ws.Range("D263").Value = 9
$ws.Range("C273").Value = 2.339181286549707
$ws.Range("D273").Value = 4
$ws.Range("C274").Value = 4.929577464788732
$ws.Range("D274").Value = 7
$ws.Range("C287").Value = 5.263157894736842
$ws.Range("D287").Value = 7
$ws.Range("C298").Value = 1.360544217687075
$ws.Range("D298").Value = 2
$ws.Range("C317").Value = 0.8771929824561403
$ws.Range("D317").Value = 1
$ws.Range("E317").Value = 114
$ws.Range("C324").Value = 2.030456852791878
$ws.Range("D324").Value = 4
$ws.Range("C330").Value = 1.063829787234043
$ws.Range("D330").Value = 1
$ws.Range("C331").Value = 2.162162162162162
$ws.Range("D331").Value = 4
$ws.Range("C334").Value = 1.612903225806452
$ws.Range("D334").Value = 1
$ws.Range("C377").Value = 2.830188679245283
$ws.Range("D377").Value = 3
$ws.Range("C384").Value = 3.64963503649635
$ws.Range("D384").Value = 5
$ws.Range("C421").Value = 35.59322033898305
$ws.Range("D421").Value = 21
$ws.Range("E421").Value = 59
$ws.Range("C454").Value = 50
$ws.Range("D454").Value = 16
$ws.Range("C484").Value = 58.62068965517241
$ws.Range("D484").Value = 17
$ws.Range("C512").Value = 17.39130434782609
$ws.Range("D512").Value = 8
$ws.Range("E512").Value = 46
$ws.Range("C517").Value = 27.53623188405797
$ws.Range("D517").Value = 19
$ws.Range("C525").Value = 23.45679012345679
$ws.Range("D525").Value = 19
$ws.Range("C561").Value = 28.57142857142857
$ws.Range("D561").Value = 12
$ws.Range("C584").Value = 2.564102564102564
$ws.Range("D584").Value = 3
$ws.Range("C589").Value = 2.72108843537415
$ws.Range("D589").Value = 4
$ws.Range("C615").Value = 2.512562814070352
$ws.Range("D615").Value = 5
$ws.Range("C635").Value = 3.260869565217391
$ws.Range("D635").Value = 3
$ws.Range("C662").Value = 3.355704697986577
$ws.Range("D662").Value = 5
$ws.Range("C681").Value = 3.418803418803419
$ws.Range("D681").Value = 4
$ws.Range("C686").Value = 2.72108843537415
$ws.Range("D686").Value = 4
$ws.Range("C712").Value = 3.015075376884422
$ws.Range("D712").Value = 6
$ws.Range("C718").Value = 2.127659574468085
$ws.Range("D718").Value = 2
$ws.Range("C719").Value = 2.162162162162162
$ws.Range("D719").Value = 4
$ws.Range("C722").Value = 1.587301587301587
$ws.Range("D722").Value = 1
$ws.Range("C732").Value = 4.301075268817205
$ws.Range("D732").Value = 4
$ws.Range("C759").Value = 4.026845637583892
$ws.Range("D759").Value = 6
$ws.Range("C765").Value = 2.830188679245283
$ws.Range("D765").Value = 3
$ws.Range("C772").Value = 4.347826086956522
$ws.Range("D772").Value = 6
$ws.Range("C802").Value = 0
$ws.Range("D802").Value = 0
$ws.Range("C809").Value = 1.507537688442211
$ws.Range("D809").Value = 3

# régions (Worksheets.Item(2))
$ws = $wb.Worksheets.Item(2)
$ws.Range("D2").Value = 7.59
$ws.Range("E2").Value = 62
$ws.Range("D4").Value = 50.36
$ws.Range("E4").Value = 141
$ws.Range("D5").Value = 2.39
$ws.Range("E5").Value = 21
$ws.Range("D6").Value = 55.34
$ws.Range("E6").Value = 228
$ws.Range("D7").Value = 5.94
$ws.Range("E7").Value = 51
$ws.Range("D8").Value = 4.3
$ws.Range("E8").Value = 38
$ws.Range("D10").Value = 3.4
$ws.Range("E10").Value = 30
$ws.Range("D11").Value = 9.470000000000001
$ws.Range("E11").Value = 34
$ws.Range("D14").Value = 0.67
$ws.Range("E14").Value = 3
$ws.Range("D16").Value = 3.57
$ws.Range("E16").Value = 15
$ws.Range("D17").Value = 0.89
$ws.Range("E17").Value = 4
$ws.Range("D29").Value = 7.88
$ws.Range("E29").Value = 32
$ws.Range("F29").Value = 406
$ws.Range("D30").Value = 16.96
$ws.Range("E30").Value = 38
$ws.Range("F30").Value = 224
$ws.Range("D38").Value = 6.84
$ws.Range("E38").Value = 44
$ws.Range("D47").Value = 5.67
$ws.Range("E47").Value = 40
$ws.Range("D52").Value = 5.29
$ws.Range("E52").Value = 41
$ws.Range("D53").Value = 1.97
$ws.Range("E53").Value = 17
$ws.Range("D55").Value = 1.52
$ws.Range("E55").Value = 13
$ws.Range("D56").Value = 7.73
$ws.Range("E56").Value = 34
$ws.Range("D74").Value = 6.02
$ws.Range("E74").Value = 52
$ws.Range("D76").Value = 38.4
$ws.Range("E76").Value = 101
$ws.Range("D83").Value = 6.89
$ws.Range("E83").Value = 59
$ws.Range("D84").Value = 18.08
$ws.Range("E84").Value = 98
$ws.Range("D85").Value = 42.6
$ws.Range("E85").Value = 118
$ws.Range("F85").Value = 277
$ws.Range("D86").Value = 1.56
$ws.Range("E86").Value = 17
$ws.Range("D88").Value = 5.65
$ws.Range("E88").Value = 51
$ws.Range("D89").Value = 2.01
$ws.Range("E89").Value = 22
$ws.Range("D90").Value = 0.91
$ws.Range("E90").Value = 10
$ws.Range("D91").Value = 1.37
$ws.Range("E91").Value = 15
$ws.Range("D92").Value = 9.16
$ws.Range("E92").Value = 93
$ws.Range("F92").Value = 1015
$ws.Range("D93").Value = 24.91
$ws.Range("E93").Value = 136
$ws.Range("D95").Value = 1.5
$ws.Range("E95").Value = 20
$ws.Range("F95").Value = 1332
$ws.Range("D97").Value = 6.68
$ws.Range("E97").Value = 74
$ws.Range("F97").Value = 1107
$ws.Range("D98").Value = 1.87
$ws.Range("E98").Value = 25
$ws.Range("D99").Value = 0.9
$ws.Range("E99").Value = 12
$ws.Range("D100").Value = 1.28
$ws.Range("E100").Value = 17
$ws.Range("D101").Value = 3.89
$ws.Range("E101").Value = 27
$ws.Range("D106").Value = 5.58
$ws.Range("E106").Value = 40
$ws.Range("D107").Value = 1.98
$ws.Range("E107").Value = 16
$ws.Range("D109").Value = 1.61
$ws.Range("E109").Value = 13

# national (Worksheets.Item(3))
$ws = $wb.Worksheets.Item(3)
$ws.Range("B2").Value = 7.13
$ws.Range("C2").Value = 549
$ws.Range("B3").Value = 22.22
$ws.Range("C3").Value = 928
$ws.Range("D3").Value = 4176
$ws.Range("B4").Value = 46.13
$ws.Range("C4").Value = 1186
$ws.Range("D4").Value = 2571
$ws.Range("B5").Value = 1.34
$ws.Range("C5").Value = 131
$ws.Range("D5").Value = 9792
$ws.Range("B6").Value = 28.5
$ws.Range("C6").Value = 1515
$ws.Range("B7").Value = 5.4
$ws.Range("C7").Value = 459
$ws.Range("D7").Value = 8503
$ws.Range("B8").Value = 1.82
$ws.Range("C8").Value = 179
$ws.Range("B9").Value = 0.87
$ws.Range("C9").Value = 85
$ws.Range("B10").Value = 1.36
$ws.Range("C10").Value = 133
